# Timesheet_Group5.xlsx - updated by sravani date:17/02/2013
#
# On the "February 2013" sheet (rows 32-35, the weekly sub-total rows
# right after the "P10 / Bug Fixing and Re-testing" block), the Q:V
# "OFF day" marker columns were filled in (Q = "OFF", R:V = "LEAVE",
# reusing the same grey-fill marker style already used by the rows
# above/below), and the W column got its weekly totals filled in too.
# The active view was also scrolled/re-selected to cell W32.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "February 2013"

# --- Reuse the existing grey "marker" format (already used for the OFF/
# LEAVE cells elsewhere on the sheet, e.g. Q28) on Q32:V35 -----------------
$ws.Range("Q28").Copy()
$ws.Range("Q32:V35").PasteSpecial(-4122)   # xlPasteFormats - formatting only
$excel.CutCopyMode = $false

# --- Fill in the marker text -------------------------------------------
$ws.Range("Q32:Q35").Value = "OFF"
$ws.Range("R32:V35").Value = "LEAVE"

# --- Weekly totals in column W ------------------------------------------
$ws.Range("W32").Value = 0
$ws.Range("W33").Value = 0
$ws.Range("W34").Value = 7
$ws.Range("W35").Value = 0

# --- Restore the saved view/selection state -----------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("W32").Select()
